$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "Price" cells: written as literal text (matches the source inlineStr
# cells). Values that Excel would otherwise auto-detect as numbers (single dot,
# all digits) are temporarily given a Text format so they round-trip as strings,
# then the format is cleared again so no visible formatting change is left behind.

$ws.Range("D2").Value = "59.006.30"
$ws.Range("D3").Value = "2.970.02"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.55"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.75"
$ws.Range("D6").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.515"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").Value = "2.965.03"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.24"
$ws.Range("D11").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.44"
$ws.Range("D14").ClearFormats()
$ws.Range("D16").Value = "3.460.01"
$ws.Range("D18").Value = "2.967.55"
$ws.Range("D19").Value = "59.093.94"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "433.43"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.53"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.718"
$ws.Range("D22").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.06"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.64"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.20"
$ws.Range("D28").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.52"
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.14"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.103"
$ws.Range("D33").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.979"
$ws.Range("D35").ClearFormats()
$ws.Range("D36").Value = "0.0₃0758"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.38"
$ws.Range("D38").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "395.37"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0350"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").Value = "2.725.81"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.12"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "33.93"
$ws.Range("D48").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.02"
$ws.Range("D51").ClearFormats()

# E-column "Volume(1h)" cells (percentages kept as formatted text, not numbers)
$ws.Range("E2").Value = "  +1.41%  "
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("E6").Value = "  +1.01%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("E11").Value = "  +7.55%  "
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("E19").Value = "  +1.79%  "
$ws.Range("E20").Value = "  +2.69%  "
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("E23").Value = "  -1.75%  "
$ws.Range("E24").Value = "  -3.01%  "
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("E28").Value = "  +6.78%  "
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("E31").Value = "  -1.13%  "
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("E33").Value = "  +5.68%  "
$ws.Range("E34").Value = "  +1.73%  "
$ws.Range("E35").Value = "  +1.83%  "
$ws.Range("E36").Value = "  +5.62%  "
$ws.Range("E37").Value = "  -1.91%  "
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("E41").Value = "  +0.95%  "
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("E44").Value = "  -3.22%  "
$ws.Range("E45").Value = "  +3.00%  "
$ws.Range("E47").Value = "  -2.45%  "
$ws.Range("E48").Value = "  +14.28%  "
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("E50").Value = "  -2.38%  "
$ws.Range("E51").Value = "  -0.78%  "

# Row 35/36: PEPE and Mantle swap positions
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
